# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    4  = 10604
    6  = 962
    7  = 106
    8  = 1309
    9  = 8216
    12 = 4
    14 = 134
    15 = 3253
    17 = 324
    19 = 128
    21 = 284
    23 = 1702
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
